$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new character row: Stickman (GameData/RSCharacterStatTable)
$ws.Range("A5").Value = "Stickman"
$ws.Range("B5").Value = 150
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 350

# Move the selection cursor to match the post-edit state (F6)
$ws.Range("F6").Select()
